# "Generate Report for Handoff" - CI regenerates the localization-status
# report. A new handoff run completed for file
# a65632bd-c226-4dfb-8e2d-40917d54d4b9, refreshing its "Latest Handoff
# Datetime" on each per-language worksheet and its "Latest Handoff Date"
# on the Overview worksheet (row 6 in each sheet).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E6").Value = "2016-03-11 10:36:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E6").Value = "2016-03-11 10:36:31"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D6").Value = "2016-36-11 10:36:31"
